# [DOCGEN][WIP] merge structure: writer and full data
#
# Replaces the single "{{ writer }}" paragraph (Screen writer information
# section) with the expanded screen_writer field block - mirroring the
# director/producer sections elsewhere in the template:
#   {{ screen_writer }}
#   {{ screen_writer.name }}
#   {{ screen_writer.other_scripts }}
#   {{ screen_writer.accomplishments }}
#   {{ screen_writer.produced_movies }}
#   {{ screen_writer.companies_worked_with }}
# plus 3 additional trailing blank paragraphs.
#
# Range.InsertXML is used (rather than Find/Replace + manual paragraph
# inserts) so the exact authored OOXML - including the w:proofErr
# gramStart/gramEnd/spellStart/spellEnd spans and the
# w:lastRenderedPageBreak marker on the .name run - is reproduced byte
# for byte.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "{{ writer }}") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '{{ writer }}' paragraph"
}

# NOTE: Range.InsertXML replaces the addressed range's contents; the final
# <w:p> supplied in the payload is consumed as the new boundary of the
# replaced range rather than emitting an extra paragraph, so the payload
# below carries one more trailing <w:p/> than the number of blank
# paragraphs that should actually land in the saved document (4 supplied
# -> 3 land), keeping the final paragraph count lined up with the diff.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_writer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">{{ </w:t></w:r><w:r><w:t>screen_writer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>name</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_writer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>other_scripts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_writer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>accomplishments</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_writer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>produced_movies</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_writer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>companies_worked_with</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>}}</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)
